$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Apply updated crypto price/volume data as scraped on Tue Aug 27 04:57:53 UTC 2024
$ws.Range("D2").Value = "63.010.66"
$ws.Range("E2").Value = "  -1.57%  "
$ws.Range("D3").Value = "2.682.40"
$ws.Range("E3").Value = "  -2.19%  "
$ws.Range("E4").Value = "  +0.00%  "
$ws.Range("D5").Value = "555.43"
$ws.Range("E5").Value = "  -2.96%  "
$ws.Range("D6").Value = "158.71"
$ws.Range("E6").Value = "  -1.30%  "
$ws.Range("E7").Value = "  +0.01%  "
$ws.Range("E8").Value = "  -0.62%  "
$ws.Range("E9").Value = "  -3.52%  "
$ws.Range("E10").Value = "  -2.13%  "
$ws.Range("D11").Value = "0.369"
$ws.Range("E11").Value = "  -4.29%  "
$ws.Range("E12").Value = "  -7.28%  "
$ws.Range("D13").Value = "3.157.44"
$ws.Range("E13").Value = "  -2.14%  "
$ws.Range("D14").Value = "26.39"
$ws.Range("E14").Value = "  -1.82%  "
$ws.Range("D15").Value = "62.895.91"
$ws.Range("E16").Value = "  -2.41%  "
$ws.Range("D17").Value = "2.684.06"
$ws.Range("E17").Value = "  -2.30%  "
$ws.Range("D18").Value = "11.94"
$ws.Range("E18").Value = "  -2.16%  "
$ws.Range("E19").Value = "  -3.86%  "
$ws.Range("D20").Value = "345.63"
$ws.Range("E20").Value = "  -2.60%  "
$ws.Range("E21").Value = "  -4.91%  "
$ws.Range("E22").Value = "  +0.00%  "
$ws.Range("D23").Value = "0.511"
$ws.Range("E23").Value = "  -2.11%  "
$ws.Range("D24").Value = "63.37"
$ws.Range("E24").Value = "  -1.58%  "
$ws.Range("D25").Value = "0.169"
$ws.Range("E25").Value = "  -1.38%  "
$ws.Range("D26").Value = "1.00"
$ws.Range("E26").Value = "  +0.27%  "
$ws.Range("E27").Value = "  -2.72%  "
$ws.Range("E28").Value = "  -6.43%  "
$ws.Range("D29").Value = "7.32"
$ws.Range("E29").Value = "  +1.07%  "
$ws.Range("E30").Value = "  +2.03%  "
$ws.Range("E31").Value = "  -1.70%  "
$ws.Range("D32").Value = "165.53"
$ws.Range("E32").Value = "  +0.65%  "
$ws.Range("B33").Value = "NEARProtocol"
$ws.Range("C33").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("D33").Value = "4.91"
$ws.Range("E33").Value = "  -0.45%  "
$ws.Range("B34").Value = "ImmutableX"
$ws.Range("C34").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("D34").Value = "1.49"
$ws.Range("E34").Value = "  +0.47%  "
$ws.Range("E36").Value = "  -3.06%  "
$ws.Range("E37").Value = "  -1.47%  "
$ws.Range("D38").Value = "349.49"
$ws.Range("D39").Value = "6.38"
$ws.Range("E39").Value = "  +0.14%  "
$ws.Range("D40").Value = "0.957"
$ws.Range("E40").Value = "  -3.56%  "
$ws.Range("D41").Value = "4.01"
$ws.Range("E41").Value = "  -2.44%  "
$ws.Range("D42").Value = "38.28"
$ws.Range("E42").Value = "  -0.97%  "
$ws.Range("E43").Value = "  -3.44%  "
$ws.Range("D44").Value = "20.86"
$ws.Range("E44").Value = "  -5.18%  "
$ws.Range("D45").Value = "0.619"
$ws.Range("E45").Value = "  -1.28%  "
$ws.Range("E46").Value = "  -4.04%  "
$ws.Range("D47").Value = "0.999"
$ws.Range("E47").Value = "  +0.01%  "
$ws.Range("D48").Value = "11.04"
$ws.Range("E48").Value = "  +0.05%  "
$ws.Range("E49").Value = "  -3.25%  "
$ws.Range("B50").Value = "Aave"
$ws.Range("C50").Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$ws.Range("D50").Value = "129.27"
$ws.Range("E50").Value = "  -3.89%  "
$ws.Range("B51").Value = "VeChain"
$ws.Range("C51").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D51").Value = "0.0242"
$ws.Range("E51").Value = "  -3.76%  "
